$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Move the "Table 508" graphic frame to its new position (size unchanged).
# (Literal point values below are chosen so that PowerPoint's internal
# 1/65536-pt fixed-point rounding reproduces the exact target EMU offsets
# 6990982/3337080 rather than the naive EMU/12700 conversion, which is off
# by a fraction of a point due to floating point rounding.)
$table = $s.Shapes.Item(3)
$table.Left = 550.4710388183594
$table.Top = 262.76220703125

# Remove the "TextBox 509" shape (the "TO DO:" note) entirely.
$s.Shapes.Item(4).Delete()
